$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 873
$ws.Range("F4").Value = 69
$ws.Range("F6").Value = 476
$ws.Range("F8").Value = 1460
$ws.Range("F9").Value = 38105
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 141
$ws.Range("F12").Value = 471
$ws.Range("F13").Value = 646
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 56
$ws.Range("F16").Value = 48
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 72
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 50
$ws.Range("F25").Value = 488
$ws.Range("F26").Value = 352
$ws.Range("F27").Value = 481
$ws.Range("F28").Value = 0
$ws.Range("F29").Value = 31
$ws.Range("F30").Value = 322
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 783
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 160
$ws.Range("F35").Value = 179
$ws.Range("F37").Value = 0
$ws.Range("F38").Value = 41
$ws.Range("F39").Value = 926
$ws.Range("F40").Value = 0

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 0
$ws.Range("F5").Value = 4348
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("F11").Value = 71
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 62
$ws.Range("F16").Value = 2
$ws.Range("F17").Value = 4345
$ws.Range("F19").Value = 0

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1647
$ws.Range("F3").Value = 0

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1647
$ws.Range("F3").Value = 406
$ws.Range("F5").Value = 873
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F13").Value = 5
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = 7868
$ws.Range("F17").Value = 141
$ws.Range("F18").Value = 471
$ws.Range("F19").Value = 67
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 71
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("F28").Value = 43
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 174
$ws.Range("F31").Value = 936
$ws.Range("F32").Value = 50
$ws.Range("F33").Value = 352
$ws.Range("F34").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("F37").Value = 322
$ws.Range("F39").Value = 0
$ws.Range("F40").Value = 62
$ws.Range("F41").Value = 337
$ws.Range("F43").Value = 179
$ws.Range("F44").Value = 157
$ws.Range("F45").Value = 926
$ws.Range("F46").Value = 314
$ws.Range("F50").Value = 314
